$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Guillermo Toloza Guzman's stand-up answers for LUNES (C) / MARTES (D)
$ws.Range("C16").Value = "Corrección casos de uso"
$ws.Range("D16").Value = "Nada"
$ws.Range("C17").Value = "Nada"
$ws.Range("D17").Value = "Asistir a la reunión con los compañeros"
$ws.Range("C18").Value = "Coursera"
$ws.Range("D18").Value = "Ninguna"

# Match the saved view's active selection from the diff
$ws.Range("D19").Select()
